$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tgfb1"
$ws.Cells.Item(2, 3).Value = "Itgb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 88.72291666666666
$ws.Cells.Item(2, 8).Value = 266.16875
$ws.Cells.Item(2, 9).Value = 0.7675060578750151
$ws.Cells.Item(2, 10).Value = 0.7675060578750152
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.069456
$ws.Cells.Item(2, 14).Value = 0.208368
$ws.Cells.Item(2, 15).Value = 0.01627409789654661
$ws.Cells.Item(2, 16).Value = 0.01627409789654661
$ws.Cells.Item(2, 17).Value = 6.1623389
$ws.Cells.Item(2, 18).Value = 55.46105009999999
$ws.Cells.Item(2, 19).Value = 0.01249046872205056
$ws.Cells.Item(2, 20).Value = 0.01249046872205056

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tgfb1"
$ws.Cells.Item(3, 3).Value = "Itgb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 88.72291666666666
$ws.Cells.Item(3, 8).Value = 266.16875
$ws.Cells.Item(3, 9).Value = 0.7675060578750151
$ws.Cells.Item(3, 10).Value = 0.7675060578750152
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.6957970000000001
$ws.Cells.Item(3, 14).Value = 2.087391
$ws.Cells.Item(3, 15).Value = 0.1630308179872645
$ws.Cells.Item(3, 16).Value = 0.1630308179872644
$ws.Cells.Item(3, 17).Value = 61.73313924791668
$ws.Cells.Item(3, 18).Value = 555.59825323125
$ws.Cells.Item(3, 19).Value = 0.1251271404255445
$ws.Cells.Item(3, 20).Value = 0.1251271404255445

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tgfb1"
$ws.Cells.Item(4, 3).Value = "Itgb6"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 88.72291666666666
$ws.Cells.Item(4, 8).Value = 266.16875
$ws.Cells.Item(4, 9).Value = 0.7675060578750151
$ws.Cells.Item(4, 10).Value = 0.7675060578750152
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.502633333333333
$ws.Cells.Item(4, 14).Value = 10.5079
$ws.Cells.Item(4, 15).Value = 0.820695084116189
$ws.Cells.Item(4, 16).Value = 0.820695084116189
$ws.Cells.Item(4, 17).Value = 310.7638453472222
$ws.Cells.Item(4, 18).Value = 2796.874608125
$ws.Cells.Item(4, 19).Value = 0.6298884487274201
$ws.Cells.Item(4, 20).Value = 0.6298884487274202

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tgfb1"
$ws.Cells.Item(5, 3).Value = "Itgb6"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.91585
$ws.Cells.Item(5, 8).Value = 53.74755
$ws.Cells.Item(5, 9).Value = 0.1549827702197958
$ws.Cells.Item(5, 10).Value = 0.1549827702197958
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.069456
$ws.Cells.Item(5, 14).Value = 0.208368
$ws.Cells.Item(5, 15).Value = 0.01627409789654661
$ws.Cells.Item(5, 16).Value = 0.01627409789654661
$ws.Cells.Item(5, 17).Value = 1.2443632776
$ws.Cells.Item(5, 18).Value = 11.1992694984
$ws.Cells.Item(5, 19).Value = 0.002522204774834945
$ws.Cells.Item(5, 20).Value = 0.002522204774834945

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tgfb1"
$ws.Cells.Item(6, 3).Value = "Itgb6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.91585
$ws.Cells.Item(6, 8).Value = 53.74755
$ws.Cells.Item(6, 9).Value = 0.1549827702197958
$ws.Cells.Item(6, 10).Value = 0.1549827702197958
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.6957970000000001
$ws.Cells.Item(6, 14).Value = 2.087391
$ws.Cells.Item(6, 15).Value = 0.1630308179872645
$ws.Cells.Item(6, 16).Value = 0.1630308179872644
$ws.Cells.Item(6, 17).Value = 12.46579468245
$ws.Cells.Item(6, 18).Value = 112.19215214205
$ws.Cells.Item(6, 19).Value = 0.02526696780286556
$ws.Cells.Item(6, 20).Value = 0.02526696780286556

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tgfb1"
$ws.Cells.Item(7, 3).Value = "Itgb6"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.91585
$ws.Cells.Item(7, 8).Value = 53.74755
$ws.Cells.Item(7, 9).Value = 0.1549827702197958
$ws.Cells.Item(7, 10).Value = 0.1549827702197958
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.502633333333333
$ws.Cells.Item(7, 14).Value = 10.5079
$ws.Cells.Item(7, 15).Value = 0.820695084116189
$ws.Cells.Item(7, 16).Value = 0.820695084116189
$ws.Cells.Item(7, 17).Value = 62.75265340500001
$ws.Cells.Item(7, 18).Value = 564.773880645
$ws.Cells.Item(7, 19).Value = 0.1271935976420953
$ws.Cells.Item(7, 20).Value = 0.1271935976420953

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tgfb1"
$ws.Cells.Item(8, 3).Value = "Itgb6"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.960212333333333
$ws.Cells.Item(8, 8).Value = 26.880637
$ws.Cells.Item(8, 9).Value = 0.077511171905189
$ws.Cells.Item(8, 10).Value = 0.07751117190518901
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.069456
$ws.Cells.Item(8, 14).Value = 0.208368
$ws.Cells.Item(8, 15).Value = 0.01627409789654661
$ws.Cells.Item(8, 16).Value = 0.01627409789654661
$ws.Cells.Item(8, 17).Value = 0.622340507824
$ws.Cells.Item(8, 18).Value = 5.601064570416
$ws.Cells.Item(8, 19).Value = 0.001261424399661099
$ws.Cells.Item(8, 20).Value = 0.001261424399661099

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tgfb1"
$ws.Cells.Item(9, 3).Value = "Itgb6"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.960212333333333
$ws.Cells.Item(9, 8).Value = 26.880637
$ws.Cells.Item(9, 9).Value = 0.077511171905189
$ws.Cells.Item(9, 10).Value = 0.07751117190518901
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.6957970000000001
$ws.Cells.Item(9, 14).Value = 2.087391
$ws.Cells.Item(9, 15).Value = 0.1630308179872645
$ws.Cells.Item(9, 16).Value = 0.1630308179872644
$ws.Cells.Item(9, 17).Value = 6.234488860896334
$ws.Cells.Item(9, 18).Value = 56.110399748067
$ws.Cells.Item(9, 19).Value = 0.01263670975885443
$ws.Cells.Item(9, 20).Value = 0.01263670975885443

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tgfb1"
$ws.Cells.Item(10, 3).Value = "Itgb6"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 8.960212333333333
$ws.Cells.Item(10, 8).Value = 26.880637
$ws.Cells.Item(10, 9).Value = 0.077511171905189
$ws.Cells.Item(10, 10).Value = 0.07751117190518901
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.502633333333333
$ws.Cells.Item(10, 14).Value = 10.5079
$ws.Cells.Item(10, 15).Value = 0.820695084116189
$ws.Cells.Item(10, 16).Value = 0.820695084116189
$ws.Cells.Item(10, 17).Value = 31.38433839247778
$ws.Cells.Item(10, 18).Value = 282.4590455323
$ws.Cells.Item(10, 19).Value = 0.06361303774667347
$ws.Cells.Item(10, 20).Value = 0.06361303774667348
